$wb = $excel.ActiveWorkbook

# --- Fix the mis-typed "2050" column header (it was showing a stray number
#     left over from a formula) on every table that has a 2015/2030/2040
#     (or 2015/2015-2030/2031-2040) header row. ---

# Sheets whose column headers are single years -> the 5th column should read "2050".
$yearHeaderSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $yearHeaderSheets) {
    $ws = $wb.Worksheets.Item($name)
    # Enter as text (leading apostrophe) so "2050" isn't reinterpreted as a
    # number, then restore the original header formatting (D1's style) so
    # the cell keeps its usual bold/centered header look instead of picking
    # up a "quote prefix" style variant.
    $ws.Range("E1").Value = "'2050"
    $ws.Range("D1").Copy() | Out-Null
    $ws.Range("E1").PasteSpecial(-4122) | Out-Null
}

# "Potencia Incremental" uses period ranges ("2015-2030", "2031-2040"), so the
# last header should read "2041-2050" -- this isn't numeric-looking so it
# stays text on its own and keeps its existing style untouched.
$wsIncremental = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIncremental.Range("E1").Value = "2041-2050"

# --- Remove the "Total" rows that were added by mistake. ---

# Row 13 on the four 12-row tables.
$totalRowSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $totalRowSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(13).Delete()
}

# Row 4 on the "Custo Total" table.
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Rows.Item(4).Delete()

$excel.CutCopyMode = 0
